$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Tuan 5 - 16: hoan thanh code."  ->  "Tuan 5 - 14: hoan thanh code."
# ------------------------------------------------------------------
$pCode = $d.Paragraphs.Item(3)
$pCode.Range.Find.Execute("16", $false, $false, $false, $false, $false, `
    $true, 1, $false, "14", 2) | Out-Null

# ------------------------------------------------------------------
# 2) "Tuan 17: Kiem thu toan bo he thong, viet tai lieu kiem thu."
#    -> "Tuan 15: Kiem thu toan bo he thong, viet tai lieu kiem thu."
# ------------------------------------------------------------------
$pTest = $d.Paragraphs.Item(4)
$pTest.Range.Find.Execute("17", $false, $false, $false, $false, $false, `
    $true, 1, $false, "15", 2) | Out-Null

# ------------------------------------------------------------------
# 3) New paragraph after it: "Tuan 16: Hoan thanh mau truong hop kiem thu"
# ------------------------------------------------------------------
$r = $d.Paragraphs.Item(4).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$d.Paragraphs.Item(5).Range.Text = "Tuần 16: Hoàn thành mẫu trường hợp kiểm thử"

# ------------------------------------------------------------------
# 4) "Tuan 18 - 19: Viet tai lieu dam bao chat luong va hoan thanh"
#    -> "Tuan 17 - 18: Viet tai lieu dam bao chat luong va hoan thanh"
#    (paragraph shifted down by one, now index 6)
# ------------------------------------------------------------------
$pDoc = $d.Paragraphs.Item(6)
$pDoc.Range.Find.Execute("18", $false, $false, $false, $false, $false, `
    $true, 1, $false, "17", 2) | Out-Null
$pDoc2 = $d.Paragraphs.Item(6)
$pDoc2.Range.Find.Execute("19", $false, $false, $false, $false, $false, `
    $true, 1, $false, "18", 2) | Out-Null

# ------------------------------------------------------------------
# 5) New paragraph after it: "Con du 2 tuan tet: co the di choi hoac code."
#    carrying the _GoBack bookmark (removed from the "Dead line" paragraph).
# ------------------------------------------------------------------
$r2 = $d.Paragraphs.Item(6).Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(7)
# append a one-char placeholder so the insertion point used for the
# bookmark is never exactly on the paragraph mark (which this host
# mishandles); the placeholder is stripped right after.
$newPara.Range.Text = "Còn dư 2 tuần tết: có thể đi chơi hoặc code.X"

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$placeholderPos = $d.Paragraphs.Item(7).Range.End - 2
$d.Bookmarks.Add("_GoBack", $d.Range($placeholderPos, $placeholderPos))
$d.Range($placeholderPos, $placeholderPos + 1).Delete() | Out-Null
